$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-27 Friday" "2024-09-28 Saturday"
Replace-Text "813×7=5691" "918×7=6426"
Replace-Text "161×7=1127" "469×4=1876"
Replace-Text "558×9=5022" "184×4=736"
Replace-Text "595×8=4760" "330×8=2640"
Replace-Text "283×3=849" "666×2=1332"
Replace-Text "102×5=510" "238×2=476"
Replace-Text "579×9=5211" "824×9=7416"
Replace-Text "455×5=2275" "232×3=696"
Replace-Text "819×3=2457" "575×9=5175"
Replace-Text "784×9=7056" "205×8=1640"
Replace-Text "676×5=3380" "524×8=4192"
Replace-Text "480×9=4320" "473×9=4257"
Replace-Text "206×2=412" "742×4=2968"
Replace-Text "405×5=2025" "360×5=1800"
Replace-Text "549×3=1647" "683×8=5464"
Replace-Text "426×5=2130" "751×6=4506"
Replace-Text "805×4=3220" "369×6=2214"
Replace-Text "568×3=1704" "303×8=2424"
Replace-Text "536×2=1072" "423×5=2115"
Replace-Text "874×4=3496" "454×6=2724"
Replace-Text "554×2=1108" "924×2=1848"
Replace-Text "344×3=1032" "469×9=4221"
Replace-Text "405×3=1215" "399×2=798"
Replace-Text "976×7=6832" "913×3=2739"
Replace-Text "112×4=448" "987×7=6909"
